$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.200282878265071
$ws.Range("D2").Value = 0.1543759122274935
$ws.Range("E2").Value = 0.1546295407942502
$ws.Range("F2").Value = 1.805833490361131
$ws.Range("G2").Value = 0.002503252704902824
$ws.Range("I2").Value = 1.268184973858858
$ws.Range("J2").Value = 0.1948882914732621
$ws.Range("K2").Value = 2.238571896726114
$ws.Range("L2").Value = 0.2152460557927753
$ws.Range("O2").Value = 4.648411802060537
$ws.Range("C3").Value = 0.1957145503483275
$ws.Range("D3").Value = 0.1498356420788696
$ws.Range("E3").Value = 0.1532361930638615
$ws.Range("F3").Value = 1.823267539609333
$ws.Range("G3").Value = 0.002506347948600474
$ws.Range("I3").Value = 1.281257982336221
$ws.Range("J3").Value = 0.1951329391446563
$ws.Range("K3").Value = 2.028888705679663
$ws.Range("L3").Value = 0.2147287509102966
$ws.Range("O3").Value = 4.709726907168474
$ws.Range("C4").Value = 0.1929644642371642
$ws.Range("D4").Value = 0.1470759178108949
$ws.Range("E4").Value = 0.1524313634453591
$ws.Range("F4").Value = 1.835099578663296
$ws.Range("G4").Value = 0.002508348510002439
$ws.Range("I4").Value = 1.290093201530034
$ws.Range("J4").Value = 0.1953636461386346
$ws.Range("K4").Value = 1.899891944254477
$ws.Range("L4").Value = 0.2144868920320064
$ws.Range("O4").Value = 4.750521326293793
$ws.Range("C5").Value = 0.1918576867104207
$ws.Range("D5").Value = 0.1459584526993183
$ws.Range("E5").Value = 0.1521161834048712
$ws.Range("F5").Value = 1.84020454445249
$ws.Range("G5").Value = 0.002509188993293085
$ws.Range("I5").Value = 1.293896725304677
$ws.Range("J5").Value = 0.1954779272508596
$ws.Range("K5").Value = 1.847265272329423
$ws.Range("L5").Value = 0.214407430195493
$ws.Range("O5").Value = 4.767935719273993
$ws.Range("C6").Value = 0.1916747500265501
$ws.Range("D6").Value = 0.145773332938262
$ws.Range("E6").Value = 0.1520646221709541
$ws.Range("F6").Value = 1.841069325851194
$ws.Range("G6").Value = 0.002509330081509956
$ws.Range("I6").Value = 1.294540559497115
$ws.Range("J6").Value = 0.1954981283264772
$ws.Range("K6").Value = 1.83852316800801
$ws.Range("L6").Value = 0.2143953903971507
$ws.Range("O6").Value = 4.770875073130441
$ws.Range("C7").Value = 0.192949481408462
$ws.Range("D7").Value = 0.1470608182238777
$ws.Range("E7").Value = 0.1524270609485257
$ws.Range("F7").Value = 1.835167279154241
$ws.Range("G7").Value = 0.002508359742594029
$ws.Range("I7").Value = 1.290143675056278
$ws.Range("J7").Value = 0.1953651052825585
$ws.Range("K7").Value = 1.899182438529976
$ws.Range("L7").Value = 0.2144857429961675
$ws.Range("O7").Value = 4.750752984090283
$ws.Range("C8").Value = 0.198696387372209
$ws.Range("D8").Value = 0.152804669149873
$ws.Range("E8").Value = 0.1541386213358429
$ws.Range("F8").Value = 1.811610627873506
$ws.Range("G8").Value = 0.002504299220740995
$ws.Range("I8").Value = 1.272524678348361
$ws.Range("J8").Value = 0.1949559548169475
$ws.Range("K8").Value = 2.16632715014714
$ws.Range("L8").Value = 0.2150519897113981
$ws.Range("O8").Value = 4.668899600170747
$ws.Range("C9").Value = 0.2103977214372748
$ws.Range("D9").Value = 0.1642869720070195
$ws.Range("E9").Value = 0.157895403163014
$ws.Range("F9").Value = 1.774371262648032
$ws.Range("G9").Value = 0.002497127115342936
$ws.Range("I9").Value = 1.244395259006389
$ws.Range("J9").Value = 0.1947913228107225
$ws.Range("K9").Value = 2.688081965167271
$ws.Range("L9").Value = 0.2167618594152856
$ws.Range("O9").Value = 4.533394066844863
$ws.Range("C10").Value = 0.2192537691241796
$ws.Range("D10").Value = 0.172852298684532
$ws.Range("E10").Value = 0.160897644328557
$ws.Range("F10").Value = 1.752484223464187
$ws.Range("G10").Value = 0.002492334880521987
$ws.Range("I10").Value = 1.227654010190534
$ws.Range("J10").Value = 0.1950580932909816
$ws.Range("K10").Value = 3.069993385234568
$ws.Range("L10").Value = 0.218381599617075
$ws.Range("O10").Value = 4.449139988894615
$ws.Range("C11").Value = 0.2233380842245793
$ws.Range("D11").Value = 0.1767761836544963
$ws.Range("E11").Value = 0.1623156336159433
$ws.Range("F11").Value = 1.743718542626965
$ws.Range("G11").Value = 0.002490257337855419
$ws.Range("I11").Value = 1.220892711862071
$ws.Range("J11").Value = 0.1952634449125341
$ws.Range("K11").Value = 3.243401085020764
$ws.Range("L11").Value = 0.2191970073090701
$ws.Range("O11").Value = 4.414145305761707
$ws.Range("C12").Value = 0.2248926242513107
$ws.Range("D12").Value = 0.1782659237127149
$ws.Range("E12").Value = 0.16286006278089
$ws.Range("F12").Value = 1.740570690760578
$ws.Range("G12").Value = 0.00248948528382106
$ws.Range("I12").Value = 1.218455438563886
$ws.Range("J12").Value = 0.1953532642062896
$ws.Range("K12").Value = 3.309016357095174
$ws.Range("L12").Value = 0.2195170420698318
$ws.Range("O12").Value = 4.401374041540436
$ws.Range("C13").Value = 0.2245574766078846
$ws.Range("D13").Value = 0.1779449119038929
$ws.Range("E13").Value = 0.1627424788459919
$ws.Range("F13").Value = 1.741241004085765
$ws.Range("G13").Value = 0.002489650908331391
$ws.Range("I13").Value = 1.218974871222613
$ws.Range("J13").Value = 0.1953333841122245
$ws.Range("K13").Value = 3.294887235656802
$ws.Range("L13").Value = 0.2194476168422028
$ws.Range("O13").Value = 4.404103174100328
$ws.Range("C14").Value = 0.223465819464991
$ws.Range("D14").Value = 0.1768986688811793
$ws.Range("E14").Value = 0.1623602747399069
$ws.Range("F14").Value = 1.743456127820195
$ws.Range("G14").Value = 0.002490193527177435
$ws.Range("I14").Value = 1.220689727855188
$ws.Range("J14").Value = 0.1952705928798366
$ws.Range("K14").Value = 3.248800324825879
$ws.Range("L14").Value = 0.2192231113833643
$ws.Range("O14").Value = 4.413084970149413
$ws.Range("C15").Value = 0.2227981734797027
$ws.Range("D15").Value = 0.176258313222263
$ws.Range("E15").Value = 0.1621271347821143
$ws.Range("F15").Value = 1.744835300477391
$ws.Range("G15").Value = 0.002490527804601667
$ws.Range("I15").Value = 1.221756162847207
$ws.Range("J15").Value = 0.1952337009921337
$ws.Range("K15").Value = 3.220564083080774
$ws.Range("L15").Value = 0.2190870601293327
$ws.Range("O15").Value = 4.418649185395509
$ws.Range("C16").Value = 0.2189879604858049
$ws.Range("D16").Value = 0.1725964077936339
$ws.Range("E16").Value = 0.1608060228870727
$ws.Range("F16").Value = 1.753081068030355
$ws.Range("G16").Value = 0.002492472708909234
$ws.Range("I16").Value = 1.228113088474998
$ws.Range("J16").Value = 0.1950463614607614
$ws.Range("K16").Value = 3.058653919189908
$ws.Range("L16").Value = 0.2183298884265028
$ws.Range("O16").Value = 4.451494142612859
$ws.Range("C17").Value = 0.2166646980545721
$ws.Range("D17").Value = 0.1703569123192352
$ws.Range("E17").Value = 0.1600089138361973
$ws.Range("F17").Value = 1.758444781855431
$ws.Range("G17").Value = 0.002493692040323433
$ws.Range("I17").Value = 1.232231843745303
$ws.Range("J17").Value = 0.1949529328733703
$ws.Range("K17").Value = 2.959241260762553
$ws.Range("L17").Value = 0.2178854842189963
$ws.Range("O17").Value = 4.472497949986689
$ws.Range("C18").Value = 0.2153336636200578
$ws.Range("D18").Value = 0.1690714057941705
$ws.Range("E18").Value = 0.159555359031085
$ws.Range("F18").Value = 1.761641919290739
$ws.Range("G18").Value = 0.002494403016567917
$ws.Range("I18").Value = 1.234681247680797
$ws.Range("J18").Value = 0.1949071026332589
$ws.Range("K18").Value = 2.902031284124973
$ws.Range("L18").Value = 0.2176372734709418
$ws.Range("O18").Value = 4.48489242334449
$ws.Range("C19").Value = 0.2148839022062674
$ws.Range("D19").Value = 0.1686366036474567
$ws.Range("E19").Value = 0.1594026398246555
$ws.Range("F19").Value = 1.762743654894329
$ws.Range("G19").Value = 0.00249464540001776
$ws.Range("I19").Value = 1.235524376520971
$ws.Range("J19").Value = 0.1948929440647689
$ws.Range("K19").Value = 2.882655851245261
$ws.Range("L19").Value = 0.2175545056151833
$ws.Range("O19").Value = 4.48914281048431
$ws.Range("C20").Value = 0.2169114712630744
$ws.Range("D20").Value = 0.1705950429038836
$ws.Range("E20").Value = 0.1600932584660022
$ws.Range("F20").Value = 1.757862204392552
$ws.Range("G20").Value = 0.002493561242372522
$ws.Range("I20").Value = 1.231785071982671
$ws.Range("J20").Value = 0.1949620602751594
$ws.Range("K20").Value = 2.969827085282645
$ws.Range("L20").Value = 0.2179320263638616
$ws.Range("O20").Value = 4.470229589726955
$ws.Range("C21").Value = 0.2237862520942855
$ws.Range("D21").Value = 0.1772058721294201
$ws.Range("E21").Value = 0.1624723350588617
$ws.Range("F21").Value = 1.742800835023019
$ws.Range("G21").Value = 0.002490033749228919
$ws.Range("I21").Value = 1.22018269064872
$ws.Range("J21").Value = 0.1952887091316171
$ws.Range("K21").Value = 3.262338565813479
$ws.Range("L21").Value = 0.2192887488855675
$ws.Range("O21").Value = 4.410433749803246
$ws.Range("C22").Value = 0.2283253096273938
$ws.Range("D22").Value = 0.181548825679414
$ws.Range("E22").Value = 0.1640707033701894
$ws.Range("F22").Value = 1.733957245721115
$ws.Range("G22").Value = 0.002487813782302303
$ws.Range("I22").Value = 1.213317393440711
$ws.Range("J22").Value = 0.1955724611501921
$ws.Range("K22").Value = 3.453215958409203
$ws.Range("L22").Value = 0.2202410438491214
$ws.Range("O22").Value = 4.37415476938861
$ws.Range("C23").Value = 0.2258985553399953
$ws.Range("D23").Value = 0.1792288936192108
$ws.Range("E23").Value = 0.1632136588487576
$ws.Range("F23").Value = 1.738585653935274
$ws.Range("G23").Value = 0.002488990824122937
$ws.Range("I23").Value = 1.216915803564611
$ws.Range("J23").Value = 0.195414594875281
$ws.Range("K23").Value = 3.351369349196887
$ws.Range("L23").Value = 0.2197267976882031
$ws.Range("O23").Value = 4.393260853203515
$ws.Range("C24").Value = 0.2167998906354427
$ws.Range("D24").Value = 0.1704873778021181
$ws.Range("E24").Value = 0.1600551115729445
$ws.Range("F24").Value = 1.758125234236971
$ws.Range("G24").Value = 0.002493620345050225
$ws.Range("I24").Value = 1.23198680371479
$ws.Range("J24").Value = 0.1949579092202143
$ws.Range("K24").Value = 2.965041409162268
$ws.Range("L24").Value = 0.2179109619711639
$ws.Range("O24").Value = 4.47125412147895
$ws.Range("C25").Value = 0.207186412533531
$ws.Range("D25").Value = 0.1611577167383729
$ws.Range("E25").Value = 0.1568364318186966
$ws.Range("F25").Value = 1.783485397540211
$ws.Range("G25").Value = 0.002498983226375641
$ws.Range("I25").Value = 1.251316442593804
$ws.Range("J25").Value = 0.1947676969873839
$ws.Range("K25").Value = 2.688081965167271
$ws.Range("L25").Value = 0.2167618594152856
$ws.Range("O25").Value = 4.567368771989265
